# Executed for 7.1+Extent Report with docs
# Rotates the "Ink / Varnish" material rows (D/H columns, rows 2-8) so that
# each row's ink name + inventory item moves to the row that previously held
# the "next" ink in the press sequence.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "PANTONE Yellow U  - UV - "
$ws.Range("H2").Value = ""

$ws.Range("D3").Value = "Cyan - UV - "
$ws.Range("H3").Value = "10001837 - 9443 PRO CYAN BW8 UV - INK"

$ws.Range("D4").Value = "Pantone-1 - UV - "
$ws.Range("H4").Value = ""

$ws.Range("D5").Value = "Yellow - UV - "
$ws.Range("H5").Value = "10001305 - PROCESS YELLOW C UV"

$ws.Range("D6").Value = "Magenta - UV - "
$ws.Range("H6").Value = "10001836 - 9442 PRO MAGENTA BW5 UV - INK"

$ws.Range("D7").Value = "Black - UV - "
$ws.Range("H7").Value = "10001817 - 9409 MIXING BLACK UV - INK"

$ws.Range("D8").Value = "Black - UV - "
$ws.Range("H8").Value = "10001817 - 9409 MIXING BLACK UV - INK"
